# Login functionality regression tests
#
# Adds a "forgot password" regression-test account to the credentials
# sheet ("test"): a new row 6 with the user label, the demo Gmail
# account (turned into a mailto hyperlink, like the other email cells
# on this sheet), and the most-recently-generated throwaway password
# used while exercising the forgot-password flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "test" sheet

# --- new row of credentials -------------------------------------------------
$ws.Range("A6").Value = "forgot password user"
$ws.Range("B6").Value = "qataskdemoaccnt@gmail.com"
$ws.Range("C6").Value = "newPT_407*602"

# Vertically center the new row (set before the hyperlink is added so the
# hyperlink's own style - which carries the underlined/blue "Hyperlink"
# font - is the one that picks up the alignment tweak for B6).
$ws.Range("A6:C6").VerticalAlignment = -4108   # xlCenter

# The email cell becomes a clickable mailto link, matching the existing
# Email/Password hyperlink cells above it.
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:qataskdemoaccnt@gmail.com")

# User-label and password cells get a monospace font to set them apart.
$ws.Range("A6,C6").Font.Name = "JetBrains Mono"
$ws.Range("A6,C6").Font.Size = 9.8
$ws.Range("A6,C6").Font.Color = 5867370   # RGB(106,135,89)

# Column A now holds a longer label ("forgot password user") than before,
# so widen it to fit.
$ws.Columns.Item(1).ColumnWidth = 35.26

# Make "test" the active sheet/tab again, with the cell below the new row
# selected (where the next credential row would go).
$ws.Activate()
$ws.Range("A7").Select()
